$d = $word.ActiveDocument

# NOTE: do all plain Find/Replace substitutions BEFORE inserting the new body/summary
# paragraph XML below, since the new paragraph text contains words like "technological"
# that would collide with the narrow "tech" -> "edu" replacement if done afterwards.

# 1. Title
$d.Content.Find.Execute("Galaxies: Cosmic Tapestry of Wonders", $true, $false, $false, $false, $false, $true, 1, $false, "Exploring Chemistry - The Building Blocks of Life", 2) | Out-Null

# 2. Author name
$d.Content.Find.Execute("Katie Leonard", $true, $false, $false, $false, $false, $true, 1, $false, "Maya Singh", 2) | Out-Null

# 3. Email - local part and domain suffix
$d.Content.Find.Execute("intelligentwriter349@digitalworld", $true, $false, $false, $false, $false, $true, 1, $false, "mrssingh@sunnyvalehigh", 2) | Out-Null
$d.Content.Find.Execute("tech", $true, $false, $false, $false, $false, $true, 1, $false, "edu", 2) | Out-Null

# 4. Body paragraph: replace the whole paragraph (multiple runs, including new sentences)
$body_frag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>Chemistry is often regarded as the central science due to its profound influence on various fields, spanning biology, medicine, and materials science</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Its study unveils the fundamental principles governing the behavior of matter at the atomic and molecular levels</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> By delving into the intricate dance of chemical reactions, we uncover the secrets behind the formation and transformation of substances, unlocking the enigmas that shape our world</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Chemistry reveals the tapestry of forces holding atoms together, enabling us to unravel the secrets of bonding and molecular structures</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>From the vibrant hues of blooming flowers to the tantalizing flavors of our favorite foods, chemistry orchestrates the symphony of life</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> It underpins the intricate mechanisms of metabolism, the process by which living organisms convert energy from food into usable forms</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Understanding chemical reactions illuminates the intricate interplay between our bodies and the environment, revealing the profound impact of nutrition, pharmaceuticals, and various environmental exposures on our health and well-being</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>The discoveries of chemistry have revolutionized the way we live</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> From the transformative power of electricity to the development of innovative materials like plastics and semiconductors, chemistry has ignited a wave of technological advancements that have shaped modern society</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> By manipulating the properties of matter, chemists have synthesized an array of materials with tailored properties, paving the way for breakthroughs in industries ranging from electronics to medicine</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> The field continues to push the boundaries of knowledge, unlocking new frontiers in energy storage, sustainable technologies, and space exploration</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document>
'@
$bodyPara = $d.Paragraphs.Item(5).Range
$bodyPara.InsertXML($body_frag) | Out-Null

# 5. Summary paragraph: replace the whole paragraph (multiple runs, incl. lastRenderedPageBreak + new sentences)
$summary_frag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>Chemistry unveils the intricacies of matter at the atomic and molecular levels, providing a foundation for understanding the behavior of substances and the forces that govern their transformations</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> It elucidates the processes underlying life, revealing the mechanisms of metabolism and the impact of nutrition, pharmaceuticals, and environmental factors on our health</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Chemistry has played a pivotal role in technological advancements, leading to the </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:lastRenderedPageBreak/><w:t>development of innovative materials and transformative technologies that have shaped modern society</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Its ongoing discoveries hold the promise of addressing global challenges and shaping a sustainable future</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Chemistry stands as a testament to the power of science to unravel the mysteries of the universe and improve the human condition</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document>
'@
$summaryPara = $d.Paragraphs.Item(7).Range
$summaryPara.InsertXML($summary_frag) | Out-Null

# 6. Append a new empty paragraph at the very end of the document body
$endRange = $d.Paragraphs.Item(7).Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter() | Out-Null
